$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reorder rows 134-137, 139-140, 142/144/145 (source data re-sync) ---
# row 134 <- old row 136
$ws.Range("B134").Value = 7483306
$ws.Range("F134").Value = "Tecnico Universitario"
$ws.Range("G134").Value = "Club Atletico Libertad"
$ws.Range("H134").Value = 1
$ws.Range("I134").Value = 1
$ws.Range("J134").Value = "D"
$ws.Range("K134").Value = 1.5
$ws.Range("L134").Value = 4.333
$ws.Range("M134").Value = 5.75
$ws.Range("N134").Value = 1.533
$ws.Range("O134").Value = 4.2
$ws.Range("P134").Value = 5.5
$ws.Range("Q134").Value = -1
$ws.Range("R134").Value = 1.925
$ws.Range("S134").Value = 1.875
$ws.Range("T134").Value = 2.25
$ws.Range("U134").Value = 1.8
$ws.Range("V134").Value = 2
$ws.Range("W134").Value = -1
$ws.Range("X134").Value = 3.2
$ws.Range("Y134").Value = -1
$ws.Range("Z134").Value = -1
$ws.Range("AA134").Value = 0.875
$ws.Range("AB134").Value = -0.5
$ws.Range("AC134").Value = 0.5

# row 135 <- old row 137
$ws.Range("B135").Value = 7483188
$ws.Range("F135").Value = "Gualaceo SC"
$ws.Range("G135").Value = "Emelec"
$ws.Range("H135").Value = 0
$ws.Range("I135").Value = 2
$ws.Range("J135").Value = "A"
$ws.Range("K135").Value = 3.6
$ws.Range("L135").Value = 3.3
$ws.Range("M135").Value = 2.05
$ws.Range("N135").Value = 2.6
$ws.Range("O135").Value = 3.25
$ws.Range("P135").Value = 2.75
$ws.Range("Q135").Value = 0
$ws.Range("R135").Value = 1.8
$ws.Range("S135").Value = 2
$ws.Range("T135").Value = 2.5
$ws.Range("U135").Value = 1.975
$ws.Range("V135").Value = 1.825
$ws.Range("W135").Value = -1
$ws.Range("X135").Value = -1
$ws.Range("Y135").Value = 1.75
$ws.Range("Z135").Value = -1
$ws.Range("AA135").Value = 1
$ws.Range("AB135").Value = -1
$ws.Range("AC135").Value = 0.825

# row 136 <- old row 135
$ws.Range("B136").Value = 7482867
$ws.Range("F136").Value = "Cumbaya FC"
$ws.Range("G136").Value = "LDU Quito"
$ws.Range("H136").Value = 1
$ws.Range("I136").Value = 2
$ws.Range("J136").Value = "A"
$ws.Range("K136").Value = 5.25
$ws.Range("L136").Value = 3.75
$ws.Range("M136").Value = 1.65
$ws.Range("N136").Value = 9
$ws.Range("O136").Value = 4.5
$ws.Range("P136").Value = 1.363
$ws.Range("Q136").Value = 1.25
$ws.Range("R136").Value = 1.975
$ws.Range("S136").Value = 1.825
$ws.Range("T136").Value = 2.5
$ws.Range("U136").Value = 1.825
$ws.Range("V136").Value = 1.975
$ws.Range("W136").Value = -1
$ws.Range("X136").Value = -1
$ws.Range("Y136").Value = 0.363
$ws.Range("Z136").Value = 0.4875
$ws.Range("AA136").Value = -0.5
$ws.Range("AB136").Value = 0.825
$ws.Range("AC136").Value = -1

# row 137 <- old row 134
$ws.Range("B137").Value = 7482832
$ws.Range("F137").Value = "Barcelona Guayaquil"
$ws.Range("G137").Value = "Guayaquil City"
$ws.Range("H137").Value = 2
$ws.Range("I137").Value = 1
$ws.Range("J137").Value = "H"
$ws.Range("K137").Value = 1.363
$ws.Range("L137").Value = 5
$ws.Range("M137").Value = 7.5
$ws.Range("N137").Value = 1.444
$ws.Range("O137").Value = 4
$ws.Range("P137").Value = 8
$ws.Range("Q137").Value = -1.25
$ws.Range("R137").Value = 2.05
$ws.Range("S137").Value = 1.75
$ws.Range("T137").Value = 2.5
$ws.Range("U137").Value = 1.95
$ws.Range("V137").Value = 1.85
$ws.Range("W137").Value = 0.444
$ws.Range("X137").Value = -1
$ws.Range("Y137").Value = -1
$ws.Range("Z137").Value = -0.5
$ws.Range("AA137").Value = 0.375
$ws.Range("AB137").Value = 0.95
$ws.Range("AC137").Value = -1

# row 139 <- old row 140
$ws.Range("B139").Value = 7528859
$ws.Range("F139").Value = "Club Atletico Libertad"
$ws.Range("G139").Value = "Cumbaya FC"
$ws.Range("H139").Value = 3
$ws.Range("I139").Value = 1
$ws.Range("J139").Value = "H"
$ws.Range("K139").Value = 1.727
$ws.Range("L139").Value = 3.5
$ws.Range("M139").Value = 4.333
$ws.Range("N139").Value = 1.4
$ws.Range("O139").Value = 4.2
$ws.Range("P139").Value = 7
$ws.Range("Q139").Value = -1.25
$ws.Range("R139").Value = 2
$ws.Range("S139").Value = 1.8
$ws.Range("T139").Value = 2.5
$ws.Range("U139").Value = 1.95
$ws.Range("V139").Value = 1.85
$ws.Range("W139").Value = 0.3999999999999999
$ws.Range("X139").Value = -1
$ws.Range("Y139").Value = -1
$ws.Range("Z139").Value = 1
$ws.Range("AA139").Value = -1
$ws.Range("AB139").Value = 0.95
$ws.Range("AC139").Value = -1

# row 140 <- old row 139
$ws.Range("B140").Value = 7528849
$ws.Range("F140").Value = "Guayaquil City"
$ws.Range("G140").Value = "Gualaceo SC"
$ws.Range("H140").Value = 0
$ws.Range("I140").Value = 2
$ws.Range("J140").Value = "A"
$ws.Range("K140").Value = 1.833
$ws.Range("L140").Value = 3.5
$ws.Range("M140").Value = 3.75
$ws.Range("N140").Value = 2.15
$ws.Range("O140").Value = 3.4
$ws.Range("P140").Value = 3
$ws.Range("Q140").Value = -0.25
$ws.Range("R140").Value = 1.825
$ws.Range("S140").Value = 1.975
$ws.Range("T140").Value = 2.5
$ws.Range("U140").Value = 1.85
$ws.Range("V140").Value = 1.95
$ws.Range("W140").Value = -1
$ws.Range("X140").Value = -1
$ws.Range("Y140").Value = 2
$ws.Range("Z140").Value = -1
$ws.Range("AA140").Value = 0.9750000000000001
$ws.Range("AB140").Value = -1
$ws.Range("AC140").Value = 0.95

# row 142 <- old row 144
$ws.Range("B142").Value = 7528848
$ws.Range("F142").Value = "Emelec"
$ws.Range("G142").Value = "Deportivo Cuenca"
$ws.Range("H142").Value = 2
$ws.Range("I142").Value = 1
$ws.Range("J142").Value = "H"
$ws.Range("K142").Value = 1.75
$ws.Range("L142").Value = 3.5
$ws.Range("M142").Value = 4.2
$ws.Range("N142").Value = 2.4
$ws.Range("O142").Value = 3.1
$ws.Range("P142").Value = 2.75
$ws.Range("Q142").Value = -0.25
$ws.Range("R142").Value = 2.05
$ws.Range("S142").Value = 1.75
$ws.Range("T142").Value = 2.25
$ws.Range("U142").Value = 1.8
$ws.Range("V142").Value = 2
$ws.Range("W142").Value = 1.4
$ws.Range("X142").Value = -1
$ws.Range("Y142").Value = -1
$ws.Range("Z142").Value = 1.05
$ws.Range("AA142").Value = -1
$ws.Range("AB142").Value = 0.8
$ws.Range("AC142").Value = -1

# row 144 <- old row 145
$ws.Range("B144").Value = 7528857
$ws.Range("F144").Value = "Universidad Catolica del Ecuador"
$ws.Range("G144").Value = "Barcelona Guayaquil"
$ws.Range("H144").Value = 0
$ws.Range("I144").Value = 1
$ws.Range("J144").Value = "A"
$ws.Range("K144").Value = 1.533
$ws.Range("L144").Value = 4
$ws.Range("M144").Value = 5.5
$ws.Range("N144").Value = 1.5
$ws.Range("O144").Value = 4.333
$ws.Range("P144").Value = 5.25
$ws.Range("Q144").Value = -1
$ws.Range("R144").Value = 1.8
$ws.Range("S144").Value = 2
$ws.Range("T144").Value = 3
$ws.Range("U144").Value = 1.975
$ws.Range("V144").Value = 1.825
$ws.Range("W144").Value = -1
$ws.Range("X144").Value = -1
$ws.Range("Y144").Value = 4.25
$ws.Range("Z144").Value = -1
$ws.Range("AA144").Value = 1
$ws.Range("AB144").Value = -1
$ws.Range("AC144").Value = 0.825

# row 145 <- old row 142
$ws.Range("B145").Value = 7528852
$ws.Range("F145").Value = "Delfin SC"
$ws.Range("G145").Value = "Tecnico Universitario"
$ws.Range("H145").Value = 2
$ws.Range("I145").Value = 2
$ws.Range("J145").Value = "D"
$ws.Range("K145").Value = 2.1
$ws.Range("L145").Value = 3.4
$ws.Range("M145").Value = 3.1
$ws.Range("N145").Value = 2.1
$ws.Range("O145").Value = 3.4
$ws.Range("P145").Value = 3.1
$ws.Range("Q145").Value = -0.25
$ws.Range("R145").Value = 1.8
$ws.Range("S145").Value = 2
$ws.Range("T145").Value = 2.25
$ws.Range("U145").Value = 1.9
$ws.Range("V145").Value = 1.9
$ws.Range("W145").Value = -1
$ws.Range("X145").Value = 2.4
$ws.Range("Y145").Value = -1
$ws.Range("Z145").Value = -0.5
$ws.Range("AA145").Value = 0.5
$ws.Range("AB145").Value = 0.8999999999999999
$ws.Range("AC145").Value = -1

# --- Update rows 185-188 with upcoming-match odds refresh (from shifted rows 188-191) ---
# row 185 update
$ws.Range("B185").Value = 7773490
$ws.Range("E185").Value = 45381.83333333334
$ws.Range("F185").Value = "Cumbaya FC"
$ws.Range("G185").Value = "El Nacional"
$ws.Range("K185").Value = 3.1
$ws.Range("L185").Value = 3.3
$ws.Range("M185").Value = 2.15
$ws.Range("N185").Value = 4.5
$ws.Range("O185").Value = 3.6
$ws.Range("P185").Value = 1.7
$ws.Range("Q185").Value = 0.75
$ws.Range("R185").Value = 1.85
$ws.Range("S185").Value = 1.95
$ws.Range("T185").Value = 2.5
$ws.Range("U185").Value = 2
$ws.Range("V185").Value = 1.8

# row 186 update
$ws.Range("B186").Value = 7773493
$ws.Range("E186").Value = 45382.625
$ws.Range("F186").Value = "Universidad Catolica del Ecuador"
$ws.Range("G186").Value = "Orense"
$ws.Range("K186").Value = 1.4
$ws.Range("L186").Value = 4.333
$ws.Range("M186").Value = 6.5
$ws.Range("N186").Value = 1.4
$ws.Range("O186").Value = 4.333
$ws.Range("P186").Value = 6.5
$ws.Range("Q186").Value = -1.25
$ws.Range("R186").Value = 1.975
$ws.Range("S186").Value = 1.825
$ws.Range("T186").Value = 2.5
$ws.Range("U186").Value = 1.85
$ws.Range("V186").Value = 1.95

# row 187 update
$ws.Range("B187").Value = 7773492
$ws.Range("E187").Value = 45382.72916666666
$ws.Range("F187").Value = "Deportivo Cuenca"
$ws.Range("G187").Value = "Delfin SC"
$ws.Range("K187").Value = 2.1
$ws.Range("L187").Value = 3.2
$ws.Range("M187").Value = 3.3
$ws.Range("N187").Value = 1.8
$ws.Range("O187").Value = 3.4
$ws.Range("P187").Value = 4.2
$ws.Range("Q187").Value = -0.5
$ws.Range("R187").Value = 1.825
$ws.Range("S187").Value = 1.975
$ws.Range("T187").Value = 2.25
$ws.Range("U187").Value = 1.9
$ws.Range("V187").Value = 1.9

# row 188 update
$ws.Range("B188").Value = 8018936
$ws.Range("E188").Value = 45382.83333333334
$ws.Range("F188").Value = "Emelec"
$ws.Range("G188").Value = "Independiente del Valle"
$ws.Range("K188").Value = 2.375
$ws.Range("L188").Value = 3.2
$ws.Range("M188").Value = 2.875
$ws.Range("N188").Value = 2.3
$ws.Range("O188").Value = 3.3
$ws.Range("P188").Value = 2.8
$ws.Range("Q188").Value = -0.25
$ws.Range("R188").Value = 2
$ws.Range("S188").Value = 1.8
$ws.Range("T188").Value = 2.5
$ws.Range("U188").Value = 1.825
$ws.Range("V188").Value = 1.975

# --- Remove now-superseded trailing rows 189-191 ---
$ws.Rows("189:191").Delete()
